# Final description label edit for P2 of Dataschema and DPE's
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C17").Value2  = "Use of contraceptive pills or injections"
$ws.Range("C24").Value2  = "History of diabetes"
$ws.Range("C35").Value2  = "Screening, skin cancer"
$ws.Range("C36").Value2  = "Screening, mammography"
$ws.Range("C37").Value2  = "Screening cervical, smear test"
$ws.Range("C59").Value2  = "Type of Cancer (ICD 10, 3 digits,e.g. C18)"
$ws.Range("C66").Value2  = "Body Mass Index at baseline"
$ws.Range("C67").Value2  = "Body Mass Index at follow-up"
$ws.Range("C68").Value2  = "Body Mass Index Standard Deviation Score at baseline (children studies)"
$ws.Range("C69").Value2  = "Body Mass Index Standard Deviation Score at follow-up (children studies)"
$ws.Range("C76").Value2  = "Body fat precent at follow-up"
$ws.Range("C77").Value2  = "Body fat precent at baseline"
$ws.Range("C94").Value2  = "Daily glycaemic load"
$ws.Range("C98").Value2  = "Intake of cakes and fine bakery products [g/d]"
$ws.Range("C99").Value2  = "Intake of fruit and vegetable juices [g/d]"
$ws.Range("C100").Value2 = "Intake of soft drinks [g/d]"
$ws.Range("C103").Value2 = "Total legumes intake [g/d]"
$ws.Range("C104").Value2 = "Total fruit intake [g/d]"

# Reflect reviewer's final selection over the edited label column
$ws.Range("C2:C109").Select()
